$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.440.51"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "2.929.17"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("E7").Value = "  -2.64%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.08%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.34%  "

$ws.Range("D14").Value = "3.390.29"
$ws.Range("E14").Value = "  -2.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").Value = "2.923.00"
$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.939"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.35%  "

$ws.Range("D18").Value = "51.366.19"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.91%  "

$ws.Range("D22").Value = "0.0₃0946"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.98%  "

$ws.Range("E26").Value = "  -4.95%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.78%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.79%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.20%  "

$ws.Range("E32").Value = "  -7.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("E35").Value = "  -3.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "34.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.60%  "

$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  -1.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.80%  "

$ws.Range("E42").Value = "  -6.02%  "

$ws.Range("E43").Value = "  -2.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("E45").Value = "  -6.24%  "

$ws.Range("E46").Value = "  -5.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.269"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.84%  "

$ws.Range("D48").Value = "2.019.19"
$ws.Range("E48").Value = "  -5.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.03%  "

$ws.Range("D51").Value = "3.209.86"
$ws.Range("E51").Value = "  -2.56%  "
